# Update Tab07 labour-market indicators for Gambia (row 50) and the
# downstream regional / income-group aggregate rows that include Gambia
# in their average (rows 61-97). Only columns F, G and H (employment
# rate - total / female / male, informal employment %) are affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab07")

$updates = @{
    "F50" = 81.599999999999994
    "G50" = 86.9
    "H50" = 76

    "F61" = 87.65
    "G61" = 90.985714285714295
    "H61" = 84.8642857142857

    "F62" = 81.897560975609807
    "G62" = 83.921951219512195
    "H62" = 79.914634146341498

    "F63" = 40.762790697674397
    "G63" = 40.327906976744202
    "H63" = 41.0162790697674

    "F64" = 56.359090909090902
    "G64" = 54.731818181818198
    "H64" = 57.572727272727299

    "F66" = 54.042519685039402
    "G66" = 54.4015748031496
    "H66" = 53.574015748031499

    "F68" = 83.94
    "G68" = 86.1
    "H68" = 81.915000000000006

    "F71" = 87.65
    "G71" = 90.985714285714295
    "H71" = 84.8642857142857

    "F82" = 81.5513513513514
    "G82" = 83.5324324324325
    "H82" = 79.583783783783801

    "F83" = 38.792307692307702
    "G83" = 38.729487179487201
    "H83" = 38.9

    "F84" = 91.0833333333333
    "G84" = 93.938888888888897
    "H84" = 88.605555555555597

    "F90" = 14.0432432432432
    "G90" = 13.4
    "H90" = 14.524324324324301

    "F91" = 89.040740740740802
    "G91" = 92.248148148148204
    "H91" = 86.551851851851893

    "F97" = 86.969230769230805
    "G97" = 90.838461538461601
    "H97" = 83.873076923076894
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
